# added more games, sped up simulate game logic, and drafted optimization logic
# Update transition-probability matrix values on Sheet1 with recalculated
# results from the larger simulated-game sample.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2061855670103093
$ws.Range("C2").Value = 0.5257731958762887
$ws.Range("J2").Value = 0.02405498281786942
$ws.Range("P2").Value = 0.134020618556701
$ws.Range("S2").Value = 0.1099656357388316
$ws.Range("B3").Value = 0.02395209580838323
$ws.Range("C3").Value = 0.05389221556886228
$ws.Range("J3").Value = 0.03592814371257485
$ws.Range("P3").Value = 0.688622754491018
$ws.Range("S3").Value = 0.1976047904191617
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.7272727272727273
$ws.Range("S4").Value = 0.2272727272727273
$ws.Range("B6").Value = 0.06465517241379311
$ws.Range("D6").Value = 0.01293103448275862
$ws.Range("E6").Value = 0.004310344827586207
$ws.Range("F6").Value = 0.06465517241379311
$ws.Range("J6").Value = 0.2672413793103448
$ws.Range("O6").Value = 0.01293103448275862
$ws.Range("Q6").Value = 0.1120689655172414
$ws.Range("R6").Value = 0.09051724137931035
$ws.Range("S6").Value = 0.3706896551724138
$ws.Range("B7").Value = 0.1132075471698113
$ws.Range("D7").Value = 0.0440251572327044
$ws.Range("F7").Value = 0.05660377358490566
$ws.Range("J7").Value = 0.06918238993710692
$ws.Range("O7").Value = 0.01886792452830189
$ws.Range("Q7").Value = 0.1761006289308176
$ws.Range("R7").Value = 0.05660377358490566
$ws.Range("S7").Value = 0.4654088050314465
$ws.Range("B8").Value = 0.1010989010989011
$ws.Range("D8").Value = 0.01318681318681319
$ws.Range("E8").Value = 0.002197802197802198
$ws.Range("F8").Value = 0.05494505494505494
$ws.Range("J8").Value = 0.1032967032967033
$ws.Range("O8").Value = 0.01098901098901099
$ws.Range("Q8").Value = 0.2065934065934066
$ws.Range("R8").Value = 0.1076923076923077
$ws.Range("S8").Value = 0.4
$ws.Range("B9").Value = 0.08875739644970414
$ws.Range("D9").Value = 0.01183431952662722
$ws.Range("F9").Value = 0.07100591715976332
$ws.Range("J9").Value = 0.1005917159763314
$ws.Range("O9").Value = 0.01775147928994083
$ws.Range("Q9").Value = 0.1597633136094675
$ws.Range("R9").Value = 0.1005917159763314
$ws.Range("S9").Value = 0.4497041420118343
$ws.Range("B10").Value = 0.1048
$ws.Range("D10").Value = 0.0216
$ws.Range("E10").Value = 0.0016
$ws.Range("F10").Value = 0.08
$ws.Range("J10").Value = 0.1104
$ws.Range("O10").Value = 0.0176
$ws.Range("Q10").Value = 0.1952
$ws.Range("R10").Value = 0.096
$ws.Range("S10").Value = 0.3728
$ws.Range("G11").Value = 0.1245283018867925
$ws.Range("J11").Value = 0.09056603773584905
$ws.Range("K11").Value = 0.1962264150943396
$ws.Range("L11").Value = 0.5773584905660377
$ws.Range("S11").Value = 0.01132075471698113
$ws.Range("G12").Value = 0.70625
$ws.Range("J12").Value = 0.21875
$ws.Range("K12").Value = 0.00625
$ws.Range("L12").Value = 0.03125
$ws.Range("S12").Value = 0.0375
$ws.Range("G13").Value = 0.5641025641025641
$ws.Range("J13").Value = 0.3846153846153846
$ws.Range("S13").Value = 0.05128205128205128
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.6666666666666666
$ws.Range("F15").Value = 0.03686635944700461
$ws.Range("H15").Value = 0.1612903225806452
$ws.Range("I15").Value = 0.08755760368663594
$ws.Range("J15").Value = 0.3594470046082949
$ws.Range("K15").Value = 0.04147465437788019
$ws.Range("N15").Value = 0.004608294930875576
$ws.Range("O15").Value = 0.06912442396313365
$ws.Range("S15").Value = 0.2396313364055299
$ws.Range("F16").Value = 0.01648351648351648
$ws.Range("H16").Value = 0.1648351648351648
$ws.Range("I16").Value = 0.06043956043956044
$ws.Range("J16").Value = 0.4340659340659341
$ws.Range("K16").Value = 0.09340659340659341
$ws.Range("M16").Value = 0.02747252747252747
$ws.Range("O16").Value = 0.04945054945054945
$ws.Range("S16").Value = 0.1538461538461539
$ws.Range("F17").Value = 0.03044496487119438
$ws.Range("H17").Value = 0.1920374707259953
$ws.Range("I17").Value = 0.08430913348946135
$ws.Range("J17").Value = 0.414519906323185
$ws.Range("K17").Value = 0.09133489461358314
$ws.Range("M17").Value = 0.0234192037470726
$ws.Range("O17").Value = 0.03044496487119438
$ws.Range("S17").Value = 0.1334894613583138
$ws.Range("F18").Value = 0.01408450704225352
$ws.Range("H18").Value = 0.2629107981220657
$ws.Range("I18").Value = 0.05633802816901409
$ws.Range("J18").Value = 0.3333333333333333
$ws.Range("K18").Value = 0.07981220657276995
$ws.Range("M18").Value = 0.0187793427230047
$ws.Range("O18").Value = 0.07511737089201878
$ws.Range("S18").Value = 0.1596244131455399
$ws.Range("F19").Value = 0.01857585139318885
$ws.Range("H19").Value = 0.1996904024767802
$ws.Range("I19").Value = 0.07198142414860681
$ws.Range("J19").Value = 0.3792569659442724
$ws.Range("K19").Value = 0.09674922600619196
$ws.Range("M19").Value = 0.01625386996904025
$ws.Range("N19").Value = 0.001547987616099071
$ws.Range("O19").Value = 0.07894736842105263
$ws.Range("S19").Value = 0.1369969040247678
